$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 341 (shifts existing rows 341-419 down to 342-420,
# dimension grows from A1:R419 to A1:R420).
$ws.Rows.Item(341).Insert()

# Populate the newly inserted row with the new data point.
$ws.Range("A341").Value = 3
$ws.Range("B341").Value = "Femacal de La Calera"
$ws.Range("C341").Value = "Coquimbo"
$ws.Range("D341").Value = 44785
$ws.Range("E341").Value = 5
$ws.Range("F341").Value = 100112031
$ws.Range("G341").Value = "Poroto verde"
$ws.Range("H341").Value = "Magnum"
$ws.Range("I341").Value = "Primera"
$ws.Range("J341").Value = 93
$ws.Range("K341").Value = 33000
$ws.Range("L341").Value = 34000
$ws.Range("M341").Value = 33516
$ws.Range("N341").Value = "$/malla 25 kilos"
$ws.Range("O341").Value = "Región de Arica y Parinacota"
$ws.Range("P341").Value = 1341
$ws.Range("Q341").Value = 25
$ws.Range("R341").Value = "Hortaliza"
